$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) -------------------------------------------------
# Written left-to-right so the shared-string table gets ID, NRIC,
# PROJECT ID, ENQUIRY, REPLY, ENQUIRY_DATE, REPLY_DATE as entries 0-6.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "NRIC"
$ws.Range("C1").Value = "PROJECT ID"
$ws.Range("D1").Value = "ENQUIRY"
$ws.Range("E1").Value = "REPLY"
$ws.Range("F1").Value = "ENQUIRY_DATE"
$ws.Range("G1").Value = "REPLY_DATE"

# --- Data rows ---------------------------------------------------------
# String cells are written in the precise order needed so new shared
# strings are interned as: S1234567A(7), Hello(8), T7654321B(9),
# "Hello to You!"(10).
$ws.Range("B2").Value = "S1234567A"
$ws.Range("D2").Value = "Hello"
$ws.Range("B3").Value = "T7654321B"
$ws.Range("E2").Value = "Hello to You!"
$ws.Range("D3").Value = "Hello"

# Numeric / date cells.
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 3
$ws.Range("F2").Value = 45764
$ws.Range("G2").Value = 45764.81342239583

$ws.Range("A3").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("F3").Value = 45764
$ws.Range("F3").NumberFormat = "mm-dd-yy"

# --- Column widths -------------------------------------------------
$ws.Columns("C").ColumnWidth = 9
$ws.Columns("F").ColumnWidth = 12
$ws.Columns("G").ColumnWidth = 9.833333333333334

# --- Selection ---------------------------------------------------------
$ws.Range("A1:G3").Select() | Out-Null
